$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.787.58"
$ws.Range("E2").Value = "  -0.56%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.123.23"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.98"
$ws.Range("E5").Value = "  -2.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.48"
$ws.Range("E6").Value = "  -3.32%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.120.66"
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.36"
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.466"
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("E13").Value = "  -1.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.06"
$ws.Range("E14").Value = "  -1.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.635.47"
$ws.Range("E15").Value = "  -1.11%  "
$ws.Range("E16").Value = "  +2.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.781.62"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.135.77"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.82"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "482.59"
$ws.Range("E20").Value = "  +1.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.72"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.706"
$ws.Range("E22").Value = "  -1.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.60"
$ws.Range("E23").Value = "  -5.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.67"
$ws.Range("E24").Value = "  +3.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.36"
$ws.Range("E25").Value = "  -2.58%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  -3.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.24"
$ws.Range("E28").Value = "  -3.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.95"
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("E30").Value = "  -2.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.04"
$ws.Range("E31").Value = "  +2.74%  "
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("E33").Value = "  -6.96%  "
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("E35").Value = "  -2.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.96"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.60"
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("E38").Value = "  -5.58%  "
$ws.Range("E39").Value = "  -6.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "438.04"
$ws.Range("E40").Value = "  -5.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0393"
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.26"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.865.84"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("E45").Value = "  -3.55%  "
$ws.Range("E46").Value = "  -4.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.38"
$ws.Range("E47").Value = "  -2.66%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.88"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.113"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.92"
$ws.Range("E51").Value = "  +2.13%  "
